# Re-generate Supplementary Table S5 cpu_time data (rows 47-96).
# A new benchmarking condition ("SP_narrower_range") is inserted for the SBS_set1
# data set, cpu_time values for some existing approach/seed combinations were
# re-measured, and the table below row 62 shifts down by 5 rows to make room.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 5 new rows (new SP_narrower_range block for SBS_set1, plus the
# extra SBS_set2/mSigHdp_ds_3k rows) by inserting 5 blank rows before row 62, which
# shifts the old rows 62-91 down to 67-96.
$ws.Rows("62:66").Insert()

# Re-measured cpu_time values for SBS_set1 / SignatureAnalyzer (rows 47-51) and the
# first SBS_set1 / signeR run (row 52); Data_set/Approach/Run columns are unchanged.
$ws.Cells.Item(47, 4).Value = 2573377.398
$ws.Cells.Item(48, 4).Value = 2451739.718
$ws.Cells.Item(49, 4).Value = 2529103.88
$ws.Cells.Item(50, 4).Value = 2533497.948
$ws.Cells.Item(51, 4).Value = 2522122.847
$ws.Cells.Item(52, 4).Value = 2500769.117

# Re-measured cpu_time values for SBS_set1 / SigProfilerExtractor (rows 57-61).
$ws.Cells.Item(57, 4).Value = 790904.3
$ws.Cells.Item(58, 4).Value = 793824.9
$ws.Cells.Item(59, 4).Value = 791003.05
$ws.Cells.Item(60, 4).Value = 789803.49
$ws.Cells.Item(61, 4).Value = 789518.97

# Rows 62-96: new SP_narrower_range block for SBS_set1, the shifted-down remainder of
# the table, and the new SBS_set1/SBS_set2 mSigHdp_ds_3k rows - write Data_set,
# Approach, Run and cpu_time for each row.
$rows62to96 = @(
  @("SBS_set1", "SP_narrower_range", "seed.1076753", 357324.52),
  @("SBS_set1", "SP_narrower_range", "seed.145879", 393020.23),
  @("SBS_set1", "SP_narrower_range", "seed.200437", 361345.89),
  @("SBS_set1", "SP_narrower_range", "seed.310111", 360435.51),
  @("SBS_set1", "SP_narrower_range", "seed.528401", 359717.25),
  @("SBS_set2", "mSigHdp", "seed.1076753", 24462112.846),
  @("SBS_set2", "mSigHdp", "seed.145879", 23800788.824),
  @("SBS_set2", "mSigHdp", "seed.200437", 25439862.31),
  @("SBS_set2", "mSigHdp", "seed.310111", 24911450.045),
  @("SBS_set2", "mSigHdp", "seed.528401", 25037565.107),
  @("SBS_set2", "SignatureAnalyzer", "seed.1076753", 6397205.75),
  @("SBS_set2", "SignatureAnalyzer", "seed.145879", 5310002.954),
  @("SBS_set2", "SignatureAnalyzer", "seed.200437", 5603919.865),
  @("SBS_set2", "SignatureAnalyzer", "seed.310111", 6500565.178),
  @("SBS_set2", "SignatureAnalyzer", "seed.528401", 5898554.229),
  @("SBS_set2", "signeR", "seed.1076753", 1538994.096),
  @("SBS_set2", "signeR", "seed.145879", 1656785.884),
  @("SBS_set2", "signeR", "seed.200437", 1569906.974),
  @("SBS_set2", "signeR", "seed.310111", 1390289.565),
  @("SBS_set2", "signeR", "seed.528401", 1442116.337),
  @("SBS_set2", "SigProfilerExtractor", "seed.1076753", 2102332.48),
  @("SBS_set2", "SigProfilerExtractor", "seed.145879", 2102739.36),
  @("SBS_set2", "SigProfilerExtractor", "seed.200437", 2107598.81),
  @("SBS_set2", "SigProfilerExtractor", "seed.310111", 2086986.95),
  @("SBS_set2", "SigProfilerExtractor", "seed.528401", 2113732.91),
  @("SBS_set1", "mSigHdp_ds_3k", "seed.1076753", 1322527.022),
  @("SBS_set1", "mSigHdp_ds_3k", "seed.145879", 1319799.728),
  @("SBS_set1", "mSigHdp_ds_3k", "seed.200437", 1329515.059),
  @("SBS_set1", "mSigHdp_ds_3k", "seed.310111", 1320033.22),
  @("SBS_set1", "mSigHdp_ds_3k", "seed.528401", 1323129.048),
  @("SBS_set2", "mSigHdp_ds_3k", "seed.1076753", 2989747.226),
  @("SBS_set2", "mSigHdp_ds_3k", "seed.145879", 3014169.645),
  @("SBS_set2", "mSigHdp_ds_3k", "seed.200437", 3007256.781),
  @("SBS_set2", "mSigHdp_ds_3k", "seed.310111", 3015428.688),
  @("SBS_set2", "mSigHdp_ds_3k", "seed.528401", 3015978.896)
)

$r = 62
foreach ($row in $rows62to96) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
